$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.989.84'
$ws.Range("E2").Value = '  +1.06%  '

$ws.Range("D3").Value = '1.640.60'
$ws.Range("E3").Value = '  +0.45%  '

$ws.Range("E4").Value = '  -0.02%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '212.80'
$cell.Style = "Normal"

$ws.Range("E6").Value = '  +0.40%  '

$ws.Range("E7").Value = '  -0.04%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '23.51'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +1.33%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.258'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -2.03%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.0881'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +2.20%  '

$ws.Range("D12").Value = '1.874.30'
$ws.Range("E12").Value = '  +0.47%  '

$ws.Range("D13").Value = '1.642.38'
$ws.Range("E13").Value = '  +0.46%  '

$ws.Range("E14").Value = '  +3.85%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '4.10'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +1.41%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '65.84'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +1.02%  '

$ws.Range("D17").Value = '27.990.33'
$ws.Range("E17").Value = '  +1.16%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '235.45'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +2.60%  '

$ws.Range("D19").Value = '0.0₃0724'
$ws.Range("E19").Value = '  +0.71%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '7.61'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +0.65%  '

$ws.Range("E21").Value = '  +0.03%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '10.62'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -0.44%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '4.37'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +0.63%  '

$ws.Range("E24").Value = '  -1.72%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '151.77'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +1.84%  '

$ws.Range("E26").Value = '  +1.33%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '15.68'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +0.72%  '

$ws.Range("E28").Value = '  +0.15%  '

$ws.Range("E29").Value = '  +0.01%  '

$ws.Range("E30").Value = '  +0.55%  '

$ws.Range("E32").Value = '  +2.11%  '

$ws.Range("E33").Value = '  +1.34%  '

$ws.Range("D34").Value = '1.419.49'
$ws.Range("E34").Value = '  -3.55%  '

$ws.Range("E35").Value = '  +2.28%  '

$ws.Range("E36").Value = '  +1.40%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '0.0169'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +1.50%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.882'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +0.70%  '

$ws.Range("E39").Value = '  -0.34%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.903'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -3.62%  '

$ws.Range("E41").Value = '  +1.14%  '

$ws.Range("E42").Value = '  -0.03%  '

$ws.Range("E43").Value = '  +6.71%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '66.65'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -1.66%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '5.53'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +3.24%  '

$ws.Range("E46").Value = '  -0.04%  '

$ws.Range("D47").Value = '1.782.79'
$ws.Range("E47").Value = '  +0.59%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '87.78'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +0.17%  '

$ws.Range("E49").Value = '  +0.95%  '

$ws.Range("E50").Value = '  +0.26%  '

$ws.Range("E51").Value = '  -1.33%  '

